# This script reproduces a scheduled market-data refresh for the Leve
# profit tracker: currentAveragePrice(NQ/HQ), LevePrice(NQ/HQ) and the
# derived LeveProfit(NQ/HQ) columns (H, I, J, K, L, M, N) are refreshed
# with newly polled Universalis price data on a handful of rows across
# several job sheets.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 15
$ws.Cells.Item(15,8).Value2 = 668.05554
$ws.Cells.Item(15,9).Value2 = 668.05554
$ws.Cells.Item(15,11).Value2 = 2004.16662
$ws.Cells.Item(15,13).Value2 = -1835.16662
# Row 17
$ws.Cells.Item(17,8).Value2 = 1118.6666
$ws.Cells.Item(17,10).Value2 = 1146
$ws.Cells.Item(17,12).Value2 = 3438
$ws.Cells.Item(17,14).Value2 = -3774
# Row 41
$ws.Cells.Item(41,8).Value2 = 912.6667
$ws.Cells.Item(41,9).Value2 = 175
$ws.Cells.Item(41,11).Value2 = 175
$ws.Cells.Item(41,13).Value2 = 265
# Row 53
$ws.Cells.Item(53,8).Value2 = 1029.6
$ws.Cells.Item(53,9).Value2 = 1032.6666
$ws.Cells.Item(53,10).Value2 = 1002
$ws.Cells.Item(53,11).Value2 = 1032.6666
$ws.Cells.Item(53,12).Value2 = 1002
$ws.Cells.Item(53,13).Value2 = -395.6666
$ws.Cells.Item(53,14).Value2 = -2276
# Row 86
$ws.Cells.Item(86,8).Value2 = 1300
$ws.Cells.Item(86,9).Value2 = 600
$ws.Cells.Item(86,11).Value2 = 600
$ws.Cells.Item(86,13).Value2 = 523
# Row 89
$ws.Cells.Item(89,8).Value2 = 1300
$ws.Cells.Item(89,9).Value2 = 600
$ws.Cells.Item(89,11).Value2 = 3000
$ws.Cells.Item(89,13).Value2 = 2616
# Row 106
$ws.Cells.Item(106,8).Value2 = 2788.125
$ws.Cells.Item(106,9).Value2 = 2383.1667
$ws.Cells.Item(106,11).Value2 = 2383.1667
$ws.Cells.Item(106,13).Value2 = -1752.1667
# Row 115
$ws.Cells.Item(115,8).Value2 = 4200
$ws.Cells.Item(115,10).Value2 = 4200
$ws.Cells.Item(115,12).Value2 = 12600
$ws.Cells.Item(115,14).Value2 = -15734
# Row 125
$ws.Cells.Item(125,8).Value2 = 2142.5
$ws.Cells.Item(125,10).Value2 = 2142.5
$ws.Cells.Item(125,12).Value2 = 19282.5
$ws.Cells.Item(125,14).Value2 = -24202.5
# Row 138
$ws.Cells.Item(138,8).Value2 = 8749
$ws.Cells.Item(138,10).Value2 = 8749
$ws.Cells.Item(138,12).Value2 = 26247
$ws.Cells.Item(138,14).Value2 = -36527

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Cells.Item(45,8).Value2 = 1471.1428
$ws.Cells.Item(45,9).Value2 = 1587
$ws.Cells.Item(45,10).Value2 = 1316.6666
$ws.Cells.Item(45,11).Value2 = 1587
$ws.Cells.Item(45,12).Value2 = 1316.6666
$ws.Cells.Item(45,13).Value2 = -1210
$ws.Cells.Item(45,14).Value2 = -2070.6666
# Row 110
$ws.Cells.Item(110,8).Value2 = 868.8
$ws.Cells.Item(110,9).Value2 = 787.25
$ws.Cells.Item(110,11).Value2 = 787.25
$ws.Cells.Item(110,13).Value2 = 1257.75
# Row 122
$ws.Cells.Item(122,8).Value2 = 35715140
$ws.Cells.Item(122,9).Value2 = 41667500
$ws.Cells.Item(122,10).Value2 = 995
$ws.Cells.Item(122,11).Value2 = 125002500
$ws.Cells.Item(122,12).Value2 = 2985
$ws.Cells.Item(122,13).Value2 = -125000050
$ws.Cells.Item(122,14).Value2 = -7885

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Cells.Item(94,8).Value2 = 1919.6666
$ws.Cells.Item(94,9).Value2 = 1038.1428
$ws.Cells.Item(94,10).Value2 = 5005
$ws.Cells.Item(94,11).Value2 = 1038.1428
$ws.Cells.Item(94,12).Value2 = 5005
$ws.Cells.Item(94,13).Value2 = -587.1428000000001
$ws.Cells.Item(94,14).Value2 = -5907
# Row 107
$ws.Cells.Item(107,8).Value2 = 1247.7142
$ws.Cells.Item(107,9).Value2 = 1206.8
$ws.Cells.Item(107,10).Value2 = 1350
$ws.Cells.Item(107,11).Value2 = 1206.8
$ws.Cells.Item(107,12).Value2 = 1350
$ws.Cells.Item(107,13).Value2 = 713.2
$ws.Cells.Item(107,14).Value2 = -5190

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Cells.Item(16,8).Value2 = 3712.5
$ws.Cells.Item(16,9).Value2 = 1773.3334
$ws.Cells.Item(16,10).Value2 = 4876
$ws.Cells.Item(16,11).Value2 = 1773.3334
$ws.Cells.Item(16,12).Value2 = 4876
$ws.Cells.Item(16,13).Value2 = -1486.3334
$ws.Cells.Item(16,14).Value2 = -5450
# Row 102
$ws.Cells.Item(102,8).Value2 = 27747.5
$ws.Cells.Item(102,10).Value2 = 27747.5
$ws.Cells.Item(102,12).Value2 = 27747.5
$ws.Cells.Item(102,14).Value2 = -32615.5
# Row 113
$ws.Cells.Item(113,8).Value2 = 3712.5
$ws.Cells.Item(113,9).Value2 = 1773.3334
$ws.Cells.Item(113,10).Value2 = 4876
$ws.Cells.Item(113,11).Value2 = 1773.3334
$ws.Cells.Item(113,12).Value2 = 4876
$ws.Cells.Item(113,13).Value2 = 396.6666
$ws.Cells.Item(113,14).Value2 = -9216
# Row 132
$ws.Cells.Item(132,8).Value2 = 6131.7
$ws.Cells.Item(132,10).Value2 = 12000
$ws.Cells.Item(132,12).Value2 = 36000
$ws.Cells.Item(132,14).Value2 = -41060
# Row 141
$ws.Cells.Item(141,8).Value2 = 1888887.4
$ws.Cells.Item(141,10).Value2 = 1888887.4
$ws.Cells.Item(141,12).Value2 = 1888887.4
$ws.Cells.Item(141,14).Value2 = -1899247.4

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 93
$ws.Cells.Item(93,8).Value2 = 30251
$ws.Cells.Item(93,10).Value2 = 30251
$ws.Cells.Item(93,12).Value2 = 30251
$ws.Cells.Item(93,14).Value2 = -33995
# Row 97
$ws.Cells.Item(97,8).Value2 = 309.8
$ws.Cells.Item(97,9).Value2 = 324.75
$ws.Cells.Item(97,10).Value2 = 250
$ws.Cells.Item(97,11).Value2 = 324.75
$ws.Cells.Item(97,12).Value2 = 250
$ws.Cells.Item(97,13).Value2 = 171.25
$ws.Cells.Item(97,14).Value2 = -1242
# Row 102
$ws.Cells.Item(102,8).Value2 = 598
$ws.Cells.Item(102,9).Value2 = 598
$ws.Cells.Item(102,11).Value2 = 598
$ws.Cells.Item(102,13).Value2 = 1024
# Row 122
$ws.Cells.Item(122,8).Value2 = 1516.6666
$ws.Cells.Item(122,9).Value2 = 1664.7142
$ws.Cells.Item(122,10).Value2 = 998.5
$ws.Cells.Item(122,11).Value2 = 4994.142599999999
$ws.Cells.Item(122,12).Value2 = 2995.5
$ws.Cells.Item(122,13).Value2 = -2544.142599999999
$ws.Cells.Item(122,14).Value2 = -7895.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Cells.Item(16,8).Value2 = 0
$ws.Cells.Item(16,9).Value2 = 0
$ws.Cells.Item(16,11).Value2 = 0
$ws.Cells.Item(16,13).ClearContents()
# Row 22
$ws.Cells.Item(22,8).Value2 = 1199.4445
$ws.Cells.Item(22,10).Value2 = 856.5714
$ws.Cells.Item(22,12).Value2 = 856.5714
$ws.Cells.Item(22,14).Value2 = -1446.5714
# Row 27
$ws.Cells.Item(27,8).Value2 = 1199.4445
$ws.Cells.Item(27,10).Value2 = 856.5714
$ws.Cells.Item(27,12).Value2 = 856.5714
$ws.Cells.Item(27,14).Value2 = -1070.5714
# Row 40
$ws.Cells.Item(40,8).Value2 = 5058
$ws.Cells.Item(40,9).Value2 = 2996
$ws.Cells.Item(40,11).Value2 = 2996
$ws.Cells.Item(40,13).Value2 = -2860
# Row 46
$ws.Cells.Item(46,8).Value2 = 516.6667
$ws.Cells.Item(46,10).Value2 = 575
$ws.Cells.Item(46,12).Value2 = 575
$ws.Cells.Item(46,14).Value2 = -951
# Row 55
$ws.Cells.Item(55,8).Value2 = 618.2857
$ws.Cells.Item(55,9).Value2 = 494.75
$ws.Cells.Item(55,11).Value2 = 494.75
$ws.Cells.Item(55,13).Value2 = -321.75
# Row 95
$ws.Cells.Item(95,8).Value2 = 50344
$ws.Cells.Item(95,10).Value2 = 50344
$ws.Cells.Item(95,12).Value2 = 50344
$ws.Cells.Item(95,14).Value2 = -55836
# Row 102
$ws.Cells.Item(102,8).Value2 = 90550
$ws.Cells.Item(102,10).Value2 = 90550
$ws.Cells.Item(102,12).Value2 = 90550
$ws.Cells.Item(102,14).Value2 = -97040
# Row 139
$ws.Cells.Item(139,8).Value2 = 80000
$ws.Cells.Item(139,10).Value2 = 80000
$ws.Cells.Item(139,12).Value2 = 80000
$ws.Cells.Item(139,14).Value2 = -90280

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 97
$ws.Cells.Item(97,8).Value2 = 21765.25
$ws.Cells.Item(97,10).Value2 = 21765.25
$ws.Cells.Item(97,12).Value2 = 21765.25
$ws.Cells.Item(97,14).Value2 = -23747.25
# Row 102
$ws.Cells.Item(102,8).Value2 = 0
$ws.Cells.Item(102,10).Value2 = 0
$ws.Cells.Item(102,12).Value2 = 0
$ws.Cells.Item(102,14).ClearContents()

